$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AE1) onto the new
# header cells so they pick up the same shared style (bold, thin border,
# centered) instead of creating a brand new style entry.
$ws.Range("AE1").Copy()
$ws.Range("AF1:AH1").PasteSpecial(-4122)

$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"
